$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing data block (rows 781:811) down by two rows so the
# previously-last two weekly entries make room for two brand-new ones
# inserted at the top of this date group.
$ws.Rows("781:782").Insert()

# Row 781 - brand new entry (Red Globe, Primera, Provincia del Elquí)
$ws.Cells.Item(781,1).Value  = 5
$ws.Cells.Item(781,2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(781,3).Value  = "Maule"
$ws.Cells.Item(781,4).Value  = 45267
$ws.Cells.Item(781,5).Value  = 7
$ws.Cells.Item(781,6).Value  = "Fruta"
$ws.Cells.Item(781,7).Value  = 100109
$ws.Cells.Item(781,8).Value  = "Uva"
$ws.Cells.Item(781,9).Value  = 100109001
$ws.Cells.Item(781,10).Value = "Uva"
$ws.Cells.Item(781,11).Value = "Red Globe"
$ws.Cells.Item(781,12).Value = "Primera"
$ws.Cells.Item(781,13).Value = 350
$ws.Cells.Item(781,14).Value = 16000
$ws.Cells.Item(781,15).Value = 16000
$ws.Cells.Item(781,16).Value = 16000
$ws.Cells.Item(781,17).Value = "`$/bandeja 8 kilos"
$ws.Cells.Item(781,18).Value = "Provincia del Elquí"
$ws.Cells.Item(781,19).Value = 2000
$ws.Cells.Item(781,20).Value = 8

# Row 782 - brand new entry (Superior Seedless, Primera, Provincia del Elquí)
$ws.Cells.Item(782,1).Value  = 5
$ws.Cells.Item(782,2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(782,3).Value  = "Maule"
$ws.Cells.Item(782,4).Value  = 45267
$ws.Cells.Item(782,5).Value  = 7
$ws.Cells.Item(782,6).Value  = "Fruta"
$ws.Cells.Item(782,7).Value  = 100109
$ws.Cells.Item(782,8).Value  = "Uva"
$ws.Cells.Item(782,9).Value  = 100109001
$ws.Cells.Item(782,10).Value = "Uva"
$ws.Cells.Item(782,11).Value = "Superior Seedless"
$ws.Cells.Item(782,12).Value = "Primera"
$ws.Cells.Item(782,13).Value = 300
$ws.Cells.Item(782,14).Value = 16000
$ws.Cells.Item(782,15).Value = 16000
$ws.Cells.Item(782,16).Value = 16000
$ws.Cells.Item(782,17).Value = "`$/bandeja 8 kilos"
$ws.Cells.Item(782,18).Value = "Provincia del Elquí"
$ws.Cells.Item(782,19).Value = 2000
$ws.Cells.Item(782,20).Value = 8
